$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Förändrad" date column (C) from 2025-07-13 (45851) to 2025-07-14 (45852)
# for every data row (rows 2 through 43).
for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45851) {
        $cell.Value2 = 45852
    }
}
